$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for "Feria Lagunitas de Puerto Montt -
# Bruselas (repollito)". It belongs chronologically/logically at the top of
# the existing detail rows (row 11), so push the current rows 11-14 down by
# one (this also grows the used range to R15) and populate the freed row 11
# with the new record.
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44781
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 100112035
$ws.Range("G11").Value = "Bruselas (repollito)"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 70
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 24000
$ws.Range("M11").Value = 24000
$ws.Range("N11").Value = "$/malla 15 kilos"
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 1600
$ws.Range("Q11").Value = 15
$ws.Range("R11").Value = "Hortaliza"
